# UCD_EngArch_Path_Electronic_ME_Modules.xlsx
# "Stage" column (E) for the Electronic/ME module rows (38-52) was storing
# the literal text "M" instead of the numeric stage. Replace it with the
# correct numeric stage values: 4 for stage-4 rows (38-44), 5 for the
# stage-5 rows (45-52). This also lets the now-unused "M" shared string
# drop out of the workbook automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 38; $r -le 44; $r++) {
    $ws.Cells.Item($r, 5).Value = 4
}
for ($r = 45; $r -le 52; $r++) {
    $ws.Cells.Item($r, 5).Value = 5
}

# Restore the selection/scroll position left by the editor onto E45.
$ws.Range("E45").Select()
